# Added ignore error types to Excel writing.
# Regenerate the randomly-generated SnippetID values (column H) on the
# "Voice Lines - main" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$snippetIds = @{
    2  = "OhHQ"
    3  = "nclm"
    4  = "D2th"
    5  = "JFUC"
    6  = "QNYz"
    7  = "0k3Z"
    8  = "PWyz"
    9  = "nitE"
    10 = "kMK7"
    11 = "kMK7"
    12 = "kMK7"
    13 = "kMK7"
    14 = "ecfI"
    15 = "JmRQ"
    16 = "4w0N"
    17 = "9GKb"
    18 = "9GKb"
    19 = "10zC"
    20 = "vdQu"
    21 = "44tN"
    22 = "vbL0"
    23 = "rclH"
    24 = "3avX"
    25 = "ak9g"
    26 = "dWrP"
    27 = "oKQv"
    28 = "K3dT"
    29 = "86MY"
    30 = "B8Oc"
    31 = "RBTG"
    32 = "a7Xy"
    33 = "zzHC"
    34 = "mKpa"
    35 = "mKpa"
    36 = "mxtg"
    37 = "mxtg"
    38 = "tBXh"
    39 = "tBXh"
    40 = "D0vq"
    41 = "D0vq"
    42 = "EwkS"
    43 = "EwkS"
    44 = "0OHN"
    45 = "0OHN"
    46 = "lrKl"
    47 = "lrKl"
    48 = "vu9c"
    49 = "CmDv"
    50 = "CmDv"
    51 = "p85p"
    52 = "p85p"
    53 = "RwXH"
    54 = "RwXH"
    55 = "RwXH"
    56 = "bSwY"
    57 = "bSwY"
    58 = "0g6w"
    59 = "X3jJ"
    60 = "D6Bb"
    61 = "Ed2a"
}

foreach ($row in $snippetIds.Keys) {
    $ws.Cells.Item($row, 8).Value = $snippetIds[$row]
}
